$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two time-series identifier headers (B1/C1) ---
$ws.Range("B1").Value = "ACIA.Flow.Inst.1Hour.0.Best-NWDM{units=CFS}"
$ws.Range("C1").Value = "ACIA.Stage.Inst.1Hour.0.Best-NWDM{units=FEET}"

# --- Column A is now hourly timestamps instead of daily dates ---
# Apply the new date/time display format to the whole working column,
# including the freshly-reserved placeholder rows below the data.
$ws.Range("A2:A16").NumberFormat = 'm/d/yy\ h:mm;@'

# First timestamp: 6/1/2021 7:00 AM
$ws.Range("A2").Value = 44348.291666666664

# Each following row rounds the previous timestamp + 1 day to the nearest hour
$ws.Range("A3").Formula = "=MROUND(A2+1,TIME(1,0,0))"
$ws.Range("A4").Formula = "=MROUND(A3+1,TIME(1,0,0))"

# Reserve empty (but formatted) rows down to row 16 for future data entry
$ws.Range("A5:A16").NumberFormat = 'm/d/yy\ h:mm;@'

# --- Column A got a touch wider to fit the new best-fit timestamp text ---
$ws.Columns.Item(1).ColumnWidth = 13.75

# --- Selection left on the newly-added placeholder rows ---
$ws.Range("A5:A16").Select()
